$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vanilla-Profile")
$ws.Range("F19").Value = "X-1P-User=(SYS_USER2)"
$ws.Rows.Item(29).Delete()
